# #CRM-41 Add acknowledge date in BB Adv Search form
# Also adds a Tracking ID column (per sharedStrings/sheet diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Tracking ID" column right after "Order ID" (new column B)
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Tracking ID"
$ws.Range("B2").Value = "{order:tracking_id}"

# Insert "Acknowledge Date" column right after "Delivery Date" (new column H)
$ws.Columns("H").Insert()
$ws.Range("H1").Value = "Acknowledge Date"
$ws.Range("H2").Value = "{order:acknowledge_date}"
